$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: FAPs | Artn | Ret | ECs
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Artn"
$ws.Range("C2").Value = "Ret"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.298902666666667
$ws.Range("H2").Value = 3.896708
$ws.Range("I2").Value = 0.7196693520699016
$ws.Range("J2").Value = 0.7196693520699017
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 2.027884333333333
$ws.Range("N2").Value = 6.083653
$ws.Range("O2").Value = 0.07019303724735296
$ws.Range("P2").Value = 0.07019303724735297
$ws.Range("Q2").Value = 2.634024368258222
$ws.Range("R2").Value = 23.706219314324
$ws.Range("S2").Value = 0.05051577763562097
$ws.Range("T2").Value = 0.05051577763562099

# Row 3: FAPs | Artn | Ret | FAPs
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Artn"
$ws.Range("C3").Value = "Ret"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.298902666666667
$ws.Range("H3").Value = 3.896708
$ws.Range("I3").Value = 0.7196693520699016
$ws.Range("J3").Value = 0.7196693520699017
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 25.63013966666666
$ws.Range("N3").Value = 76.89041899999999
$ws.Range("O3").Value = 0.887159745112283
$ws.Range("P3").Value = 0.887159745112283
$ws.Range("Q3").Value = 33.29105676007244
$ws.Range("R3").Value = 299.619510840652
$ws.Range("S3").Value = 0.6384616789474558
$ws.Range("T3").Value = 0.6384616789474559

# Row 4: FAPs | Artn | Ret | sCs
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Artn"
$ws.Range("C4").Value = "Ret"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.298902666666667
$ws.Range("H4").Value = 3.896708
$ws.Range("I4").Value = 0.7196693520699016
$ws.Range("J4").Value = 0.7196693520699017
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.232082666666667
$ws.Range("N4").Value = 3.696248
$ws.Range("O4").Value = 0.04264721764036409
$ws.Range("P4").Value = 0.04264721764036408
$ws.Range("Q4").Value = 1.600355461287111
$ws.Range("R4").Value = 14.403199151584
$ws.Range("S4").Value = 0.0306918954868249
$ws.Range("T4").Value = 0.0306918954868249

# Row 5: sCs | Artn | Ret | ECs
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Artn"
$ws.Range("C5").Value = "Ret"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.5059576666666666
$ws.Range("H5").Value = 1.517873
$ws.Range("I5").Value = 0.2803306479300983
$ws.Range("J5").Value = 0.2803306479300984
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 2.027884333333333
$ws.Range("N5").Value = 6.083653
$ws.Range("O5").Value = 0.07019303724735296
$ws.Range("P5").Value = 0.07019303724735297
$ws.Range("Q5").Value = 1.026023625563222
$ws.Range("R5").Value = 9.234212630069001
$ws.Range("S5").Value = 0.01967725961173198
$ws.Range("T5").Value = 0.01967725961173199

# Row 6: sCs | Artn | Ret | FAPs
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Artn"
$ws.Range("C6").Value = "Ret"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.5059576666666666
$ws.Range("H6").Value = 1.517873
$ws.Range("I6").Value = 0.2803306479300983
$ws.Range("J6").Value = 0.2803306479300984
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 25.63013966666666
$ws.Range("N6").Value = 76.89041899999999
$ws.Range("O6").Value = 0.887159745112283
$ws.Range("P6").Value = 0.887159745112283
$ws.Range("Q6").Value = 12.96776566208744
$ws.Range("R6").Value = 116.709890958787
$ws.Range("S6").Value = 0.2486980661648272
$ws.Range("T6").Value = 0.2486980661648272

# Row 7: sCs | Artn | Ret | sCs
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Artn"
$ws.Range("C7").Value = "Ret"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.5059576666666666
$ws.Range("H7").Value = 1.517873
$ws.Range("I7").Value = 0.2803306479300983
$ws.Range("J7").Value = 0.2803306479300984
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.232082666666667
$ws.Range("N7").Value = 3.696248
$ws.Range("O7").Value = 0.04264721764036409
$ws.Range("P7").Value = 0.04264721764036408
$ws.Range("Q7").Value = 0.6233816711671111
$ws.Range("R7").Value = 5.610435040504
$ws.Range("S7").Value = 0.01195532215353918
$ws.Range("T7").Value = 0.01195532215353918

Write-Output "done"
